$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NAME")
$ws.Range("Z2").Formula = "=1+2"
$ws.Range("Z3").Value = 42
$ws.Range("Z3").NumberFormat = "0.00"
